{"js": "// Replace the body-text run of the resume's second paragraph with the\n// new (rewritten) resume content. Each array entry becomes one <w:t> run\n// of text; consecutive entries are separated by a manual line break\n// (the same \"\\v\" / Word.InsertLocation line-break character Word itself\n// uses for <w:br/>), matching the original run/br/run/br/... structure.\nconst parts = [\n  \"Aspiring AI/ML Engineer,\",\n  \"aiming to pursue MCA in\",\n  \"AI & MLABHINAV PANDEYPhone: 9519654927, 9795626571\",\n  \"Email:abhinavpandey56393@gmail.com\",\n  \"Gender: Male\",\n  \"DOB: 24-01-2004\",\n  \"linkedin:www.linkedin.com/in/abhinav-\",\n  \"pandey-560447262\",\n  \"GitHub:https://github.com/abhinavpandey\",\n  \"010\",\n  \"ABOUT ME\",\n  \"Graduation: Dev Sanskrati vishwavidyalaya\",\n  \"Graduation (2022 - 2025): BCA (Bachelor of Computer Application)\",\n  \"Intermediate : CBSE (93%)\",\n  \"Matriculation: CBSE (88%)\",\n  \"KEY SKILLED AREAS\",\n  \"CAREER OBJECTIVE\",\n  \"Aiming to leverage my expertise in AI/ML to develop innovative olutions as\",\n  \"an AI/ML Engineer.\",\n  \"Currently working on a Gender Recognition Web App using deep learning\",\n  \"techniques.\",\n  \"I am eager to apply my skills in machine learning ,Python,and flask to create\",\n  \"impactful systems and contribute to cutting- edge technologies.Python language, MySQL, HTML, CSS, Java Script, Git\",\n  \"TECHNICALSKILLS\",\n  \"Programing Languages:Python\",\n  \"Other language: HTML,CSS,Java script\",\n  \"Frameworks/libraries:Flask,NumPy,Pandas\",\n  \"Database: Mysql\",\n  \" PROJECTS\",\n  \"Group Photo Gender Recognition Web App\",\n  \"Technologies: Flask, TensorFlow,Keras,OpenCV,MySQL\",\n  \"DESCRIPTION: A web app that recogniges gender in group photos and \",\n  \"visualizes  results.The app stores metadata and visualizes it to analyze group data.\",\n  \"Upcoming Features:Real-time detection,live video analysis.\",\n  \"OTHER FRONT -END PROJECTS\",\n  \"Created small projects to enhance my HTML,CSS,and Javascript skills.\",\n  \"         All projects are available on GitHub.\",\n  \"Acted as liaison for the university sports and annual functions, including cricket.\",\n  \"SOCIAL INTERNSHIP EXPERIENCE\",\n  \"Social internship in Kondagaon District,Bastar Region, Chhattisgarh.\",\n  \"Worked with Gayatri pariwar volunteers on awarness programs in local schools.\",\n  \"Delivered 70+ lectures on anti-addiction awarness and persnality development.\",\n  \"Assisted in organizing programs that helped local villagers including anti- substance\",\n  \"abuse campains.\",\n  \"Observed challenges and gaps in community awarness and highlighted the need for\",\n  \"more regular programs.\",\n  \"ADDITIONAL SKILLS AND INTERESTS\",\n  \"Leadership: Strong leadership skills developed through NCC and university sports.\",\n  \"Languages: Fluent in English and Hindi.\",\n  \"Interests:AI/ML,web devlopement,leadership training,cricket,and social causes.Led 160 cadets ,ensuring efficient training and fostering leadership skills within\",\n  \"the group.\",\n  \"In charge of coordinating activities ,managing logistics, and conducting training \",\n  \"sessions for cadets.Senior under officer,NCCNCC Experience\"\n];\nconst fullText = parts.join(\"\\v\");\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// The resume body (ABOUT ME / contact info / education / projects / etc.)\n// lives in the second paragraph of the document (index 1).\nconst target = paragraphs.items[1];\ntarget.clear();\nawait context.sync();\n\ntarget.insertText(fullText, Word.InsertLocation.start);\nawait context.sync();\n", "ps1": "# Replace the body-text run of the resume's second paragraph with the\n# new (rewritten) resume content. Each array entry becomes one <w:t> run\n# of text; consecutive entries are separated by a manual line break\n# (character code 11, Word's internal vertical-tab line-break, which is\n# what <w:br/> round-trips to/from), matching the original\n# run/br/run/br/... structure.\n$d = $word.ActiveDocument\n$p = $d.Paragraphs(2)\n$r = $p.Range\n$nl = [char]11\n$parts = @(\n    'Aspiring AI/ML Engineer,',\n    'aiming to pursue MCA in',\n    'AI & MLABHINAV PANDEYPhone: 9519654927, 9795626571',\n    'Email:abhinavpandey56393@gmail.com',\n    'Gender: Male',\n    'DOB: 24-01-2004',\n    'linkedin:www.linkedin.com/in/abhinav-',\n    'pandey-560447262',\n    'GitHub:https://github.com/abhinavpandey',\n    '010',\n    'ABOUT ME',\n    'Graduation: Dev Sanskrati vishwavidyalaya',\n    'Graduation (2022 - 2025): BCA (Bachelor of Computer Application)',\n    'Intermediate : CBSE (93%)',\n    'Matriculation: CBSE (88%)',\n    'KEY SKILLED AREAS',\n    'CAREER OBJECTIVE',\n    'Aiming to leverage my expertise in AI/ML to develop innovative olutions as',\n    'an AI/ML Engineer.',\n    'Currently working on a Gender Recognition Web App using deep learning',\n    'techniques.',\n    'I am eager to apply my skills in machine learning ,Python,and flask to create',\n    'impactful systems and contribute to cutting- edge technologies.Python language, MySQL, HTML, CSS, Java Script, Git',\n    'TECHNICALSKILLS',\n    'Programing Languages:Python',\n    'Other language: HTML,CSS,Java script',\n    'Frameworks/libraries:Flask,NumPy,Pandas',\n    'Database: Mysql',\n    ' PROJECTS',\n    'Group Photo Gender Recognition Web App',\n    'Technologies: Flask, TensorFlow,Keras,OpenCV,MySQL',\n    'DESCRIPTION: A web app that recogniges gender in group photos and ',\n    'visualizes  results.The app stores metadata and visualizes it to analyze group data.',\n    'Upcoming Features:Real-time detection,live video analysis.',\n    'OTHER FRONT -END PROJECTS',\n    'Created small projects to enhance my HTML,CSS,and Javascript skills.',\n    '         All projects are available on GitHub.',\n    'Acted as liaison for the university sports and annual functions, including cricket.',\n    'SOCIAL INTERNSHIP EXPERIENCE',\n    'Social internship in Kondagaon District,Bastar Region, Chhattisgarh.',\n    'Worked with Gayatri pariwar volunteers on awarness programs in local schools.',\n    'Delivered 70+ lectures on anti-addiction awarness and persnality development.',\n    'Assisted in organizing programs that helped local villagers including anti- substance',\n    'abuse campains.',\n    'Observed challenges and gaps in community awarness and highlighted the need for',\n    'more regular programs.',\n    'ADDITIONAL SKILLS AND INTERESTS',\n    'Leadership: Strong leadership skills developed through NCC and university sports.',\n    'Languages: Fluent in English and Hindi.',\n    'Interests:AI/ML,web devlopement,leadership training,cricket,and social causes.Led 160 cadets ,ensuring efficient training and fostering leadership skills within',\n    'the group.',\n    'In charge of coordinating activities ,managing logistics, and conducting training ',\n    'sessions for cadets.Senior under officer,NCCNCC Experience'\n)\n$r.Text = [string]::Join($nl, $parts)\n"}
